# Updated cryptos list on Sat Jan 20 08:10:41 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.529.70"
$ws.Range("E2").Value = "'  +0.51%  "
$ws.Range("D3").Value = "'2.477.97"
$ws.Range("E3").Value = "'  +0.55%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.10%  "
$ws.Range("D5").Value = "'313.02"
$ws.Range("E5").Value = "'  +0.51%  "
$ws.Range("D6").Value = "'92.06"
$ws.Range("E6").Value = "'  -2.42%  "
$ws.Range("D7").Value = "'0.550"
$ws.Range("E7").Value = "'  +0.61%  "
$ws.Range("E8").Value = "'  -0.27%  "
$ws.Range("E9").Value = "'  +3.50%  "
$ws.Range("D10").Value = "'32.70"
$ws.Range("E10").Value = "'  -2.20%  "
$ws.Range("E11").Value = "'  +1.57%  "
$ws.Range("E12").Value = "'  +2.02%  "
$ws.Range("B13").Value = "'Chainlink"
$ws.Range("C13").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'16.51"
$ws.Range("E13").Value = "'  +9.97%  "
$ws.Range("B14").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "'2.861.70"
$ws.Range("E14").Value = "'  +0.60%  "
$ws.Range("E15").Value = "'  -1.01%  "
$ws.Range("D16").Value = "'2.461.99"
$ws.Range("E16").Value = "'  +0.06%  "
$ws.Range("D17").Value = "'0.775"
$ws.Range("E17").Value = "'  -1.50%  "
$ws.Range("D18").Value = "'41.544.82"
$ws.Range("E18").Value = "'  +0.56%  "
$ws.Range("D19").Value = "'6.53"
$ws.Range("E19").Value = "'  +3.49%  "
$ws.Range("E20").Value = "'  +2.29%  "
$ws.Range("D21").Value = "'72.24"
$ws.Range("E21").Value = "'  +5.28%  "
$ws.Range("D22").Value = "'11.14"
$ws.Range("E22").Value = "'  -0.24%  "
$ws.Range("D23").Value = "'236.36"
$ws.Range("E23").Value = "'  -0.15%  "
$ws.Range("D24").Value = "'2.72"
$ws.Range("E24").Value = "'  -0.72%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("D26").Value = "'1.90"
$ws.Range("E26").Value = "'  +0.20%  "
$ws.Range("E27").Value = "'  +3.79%  "
$ws.Range("D28").Value = "'2.22"
$ws.Range("E28").Value = "'  -0.15%  "
$ws.Range("E29").Value = "'  +0.54%  "
$ws.Range("D30").Value = "'35.82"
$ws.Range("E30").Value = "'  -1.38%  "
$ws.Range("D31").Value = "'157.43"
$ws.Range("E31").Value = "'  +3.86%  "
$ws.Range("D32").Value = "'5.47"
$ws.Range("E32").Value = "'  -0.33%  "
$ws.Range("E33").Value = "'  -0.45%  "
$ws.Range("E34").Value = "'  +1.82%  "
$ws.Range("D35").Value = "'17.37"
$ws.Range("E35").Value = "'  +1.36%  "
$ws.Range("E36").Value = "'  -8.95%  "
$ws.Range("D37").Value = "'0.105"
$ws.Range("E37").Value = "'  +2.81%  "
$ws.Range("D38").Value = "'2.89"
$ws.Range("E38").Value = "'  -4.70%  "
$ws.Range("D39").Value = "'1.82"
$ws.Range("E39").Value = "'  -2.58%  "
$ws.Range("E40").Value = "'  +0.08%  "
$ws.Range("D41").Value = "'4.06"
$ws.Range("E41").Value = "'  -4.12%  "
$ws.Range("D43").Value = "'1.960.65"
$ws.Range("E43").Value = "'  -1.07%  "
$ws.Range("D44").Value = "'18.98"
$ws.Range("E44").Value = "'  -3.58%  "
$ws.Range("E45").Value = "'  -0.23%  "
$ws.Range("D46").Value = "'2.95"
$ws.Range("E46").Value = "'  -2.29%  "
$ws.Range("D47").Value = "'8.95"
$ws.Range("E47").Value = "'  +2.59%  "
$ws.Range("D48").Value = "'2.720.86"
$ws.Range("D49").Value = "'97.63"
$ws.Range("E49").Value = "'  +1.25%  "
$ws.Range("D50").Value = "'67.90"
$ws.Range("E50").Value = "'  -2.14%  "
$ws.Range("B51").Value = "'BitcoinSV"
$ws.Range("C51").Value = "'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'72.38"
$ws.Range("E51").Value = "'  -2.86%  "
